$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: "San Diego" campus was renamed to "UCSD" ---
$ws.Range("B7:B11").Value = "UCSD"

# --- Widen column C to fit the new/longer label ---
$ws.Columns.Item(3).ColumnWidth = 21.75

# --- Re-apply a sort over the data (by campus, then by year) ---
$dataRange = $ws.Range("A2:K11")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B11")) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("F2:F11")) | Out-Null
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# --- Update the active selection / scroll position on the sheet ---
$ws.Range("J10").Select() | Out-Null
